$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.343956666666666
$ws.Range("H2").Value = 13.03187
$ws.Range("I2").Value = 0.2551833209483726
$ws.Range("J2").Value = 0.2551833209483726
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.563107
$ws.Range("N2").Value = 7.689321
$ws.Range("O2").Value = 0.09258576031692413
$ws.Range("P2").Value = 0.09258576031692412
$ws.Range("Q2").Value = 11.13402574003
$ws.Range("R2").Value = 100.20623166027
$ws.Range("S2").Value = 0.02362634179020275
$ws.Range("T2").Value = 0.02362634179020275

$ws.Range("G3").Value = 4.343956666666666
$ws.Range("H3").Value = 13.03187
$ws.Range("I3").Value = 0.2551833209483726
$ws.Range("J3").Value = 0.2551833209483726
$ws.Range("O3").Value = 0.160803024221502
$ws.Range("P3").Value = 0.160803024221502
$ws.Range("Q3").Value = 19.33758500905889
$ws.Range("R3").Value = 174.03826508153
$ws.Range("S3").Value = 0.04103424973938449
$ws.Range("T3").Value = 0.04103424973938449

$ws.Range("G4").Value = 4.343956666666666
$ws.Range("H4").Value = 13.03187
$ws.Range("I4").Value = 0.2551833209483726
$ws.Range("J4").Value = 0.2551833209483726
$ws.Range("M4").Value = 0.7887020000000001
$ws.Range("N4").Value = 2.366106
$ws.Range("O4").Value = 0.02848986575023154
$ws.Range("P4").Value = 0.02848986575023154
$ws.Range("Q4").Value = 3.426087310913334
$ws.Range("R4").Value = 30.83478579822
$ws.Range("S4").Value = 0.007270138555517384
$ws.Range("T4").Value = 0.007270138555517384

$ws.Range("G5").Value = 4.343956666666666
$ws.Range("H5").Value = 13.03187
$ws.Range("I5").Value = 0.2551833209483726
$ws.Range("J5").Value = 0.2551833209483726
$ws.Range("M5").Value = 19.880183
$ws.Range("N5").Value = 59.640549
$ws.Range("O5").Value = 0.7181213497113423
$ws.Range("P5").Value = 0.7181213497113423
$ws.Range("Q5").Value = 86.35865347740332
$ws.Range("R5").Value = 777.22788129663
$ws.Range("S5").Value = 0.183252590863268
$ws.Range("T5").Value = 0.183252590863268

$ws.Range("I6").Value = 0.3783451411951115
$ws.Range("J6").Value = 0.3783451411951115
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.563107
$ws.Range("N6").Value = 7.689321
$ws.Range("O6").Value = 0.09258576031692413
$ws.Range("P6").Value = 0.09258576031692412
$ws.Range("Q6").Value = 16.507758128651
$ws.Range("R6").Value = 148.569823157859
$ws.Range("S6").Value = 0.03502937255976341
$ws.Range("T6").Value = 0.03502937255976341

$ws.Range("I7").Value = 0.3783451411951115
$ws.Range("J7").Value = 0.3783451411951115
$ws.Range("O7").Value = 0.160803024221502
$ws.Range("P7").Value = 0.160803024221502
$ws.Range("S7").Value = 0.06083904290368512
$ws.Range("T7").Value = 0.06083904290368512

$ws.Range("I8").Value = 0.3783451411951115
$ws.Range("J8").Value = 0.3783451411951115
$ws.Range("M8").Value = 0.7887020000000001
$ws.Range("N8").Value = 2.366106
$ws.Range("O8").Value = 0.02848986575023154
$ws.Range("P8").Value = 0.02848986575023154
$ws.Range("Q8").Value = 5.079656000152668
$ws.Range("R8").Value = 45.716904001374
$ws.Range("S8").Value = 0.01077900227990112
$ws.Range("T8").Value = 0.01077900227990112

$ws.Range("I9").Value = 0.3783451411951115
$ws.Range("J9").Value = 0.3783451411951115
$ws.Range("M9").Value = 19.880183
$ws.Range("N9").Value = 59.640549
$ws.Range("O9").Value = 0.7181213497113423
$ws.Range("P9").Value = 0.7181213497113423
$ws.Range("Q9").Value = 128.0388421229857
$ws.Range("R9").Value = 1152.349579106871
$ws.Range("S9").Value = 0.2716977234517619
$ws.Range("T9").Value = 0.2716977234517619

$ws.Range("G10").Value = 0.3495363333333334
$ws.Range("H10").Value = 1.048609
$ws.Range("I10").Value = 0.02053331770470026
$ws.Range("J10").Value = 0.02053331770470026
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.563107
$ws.Range("N10").Value = 7.689321
$ws.Range("O10").Value = 0.09258576031692413
$ws.Range("P10").Value = 0.09258576031692412
$ws.Range("Q10").Value = 0.8958990227210002
$ws.Range("R10").Value = 8.063091204489
$ws.Range("S10").Value = 0.001901092831518633
$ws.Range("T10").Value = 0.001901092831518633

$ws.Range("G11").Value = 0.3495363333333334
$ws.Range("H11").Value = 1.048609
$ws.Range("I11").Value = 0.02053331770470026
$ws.Range("J11").Value = 0.02053331770470026
$ws.Range("O11").Value = 0.160803024221502
$ws.Range("P11").Value = 0.160803024221502
$ws.Range("Q11").Value = 1.555998155196778
$ws.Range("R11").Value = 14.003983396771
$ws.Range("S11").Value = 0.003301819584216712
$ws.Range("T11").Value = 0.003301819584216712

$ws.Range("G12").Value = 0.3495363333333334
$ws.Range("H12").Value = 1.048609
$ws.Range("I12").Value = 0.02053331770470026
$ws.Range("J12").Value = 0.02053331770470026
$ws.Range("M12").Value = 0.7887020000000001
$ws.Range("N12").Value = 2.366106
$ws.Range("O12").Value = 0.02848986575023154
$ws.Range("P12").Value = 0.02848986575023154
$ws.Range("Q12").Value = 0.2756800051726668
$ws.Range("R12").Value = 2.481120046554
$ws.Range("S12").Value = 0.0005849914648137628
$ws.Range("T12").Value = 0.0005849914648137628

$ws.Range("G13").Value = 0.3495363333333334
$ws.Range("H13").Value = 1.048609
$ws.Range("I13").Value = 0.02053331770470026
$ws.Range("J13").Value = 0.02053331770470026
$ws.Range("M13").Value = 19.880183
$ws.Range("N13").Value = 59.640549
$ws.Range("O13").Value = 0.7181213497113423
$ws.Range("P13").Value = 0.7181213497113423
$ws.Range("Q13").Value = 6.948846271815667
$ws.Range("R13").Value = 62.53961644634101
$ws.Range("S13").Value = 0.01474541382415115
$ws.Range("T13").Value = 0.01474541382415115

$ws.Range("G14").Value = 5.888867
$ws.Range("H14").Value = 17.666601
$ws.Range("I14").Value = 0.3459382201518156
$ws.Range("J14").Value = 0.3459382201518156
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.563107
$ws.Range("N14").Value = 7.689321
$ws.Range("O14").Value = 0.09258576031692413
$ws.Range("P14").Value = 0.09258576031692412
$ws.Range("Q14").Value = 15.093796229769
$ws.Range("R14").Value = 135.844166067921
$ws.Range("S14").Value = 0.03202895313543933
$ws.Range("T14").Value = 0.03202895313543933

$ws.Range("G15").Value = 5.888867
$ws.Range("H15").Value = 17.666601
$ws.Range("I15").Value = 0.3459382201518156
$ws.Range("J15").Value = 0.3459382201518156
$ws.Range("O15").Value = 0.160803024221502
$ws.Range("P15").Value = 0.160803024221502
$ws.Range("Q15").Value = 26.21491763335767
$ws.Range("R15").Value = 235.934258700219
$ws.Range("S15").Value = 0.05562791199421571
$ws.Range("T15").Value = 0.05562791199421571

$ws.Range("G16").Value = 5.888867
$ws.Range("H16").Value = 17.666601
$ws.Range("I16").Value = 0.3459382201518156
$ws.Range("J16").Value = 0.3459382201518156
$ws.Range("M16").Value = 0.7887020000000001
$ws.Range("N16").Value = 2.366106
$ws.Range("O16").Value = 0.02848986575023154
$ws.Range("P16").Value = 0.02848986575023154
$ws.Range("Q16").Value = 4.644561180634001
$ws.Range("R16").Value = 41.801050625706
$ws.Range("S16").Value = 0.009855733449999269
$ws.Range("T16").Value = 0.009855733449999269

$ws.Range("G17").Value = 5.888867
$ws.Range("H17").Value = 17.666601
$ws.Range("I17").Value = 0.3459382201518156
$ws.Range("J17").Value = 0.3459382201518156
$ws.Range("M17").Value = 19.880183
$ws.Range("N17").Value = 59.640549
$ws.Range("O17").Value = 0.7181213497113423
$ws.Range("P17").Value = 0.7181213497113423
$ws.Range("Q17").Value = 117.071753622661
$ws.Range("R17").Value = 1053.645782603949
$ws.Range("S17").Value = 0.2484256215721613
$ws.Range("T17").Value = 0.2484256215721613
